$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.003", "1.000") are preserved exactly as text, matching source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.145.91"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.748.55"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "242.47"
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.5375"
$ws.Range("E7").Value = "  +3.43%  "
$ws.Range("D8").Value = "0.2840"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "0.06183"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "1.759.89"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "0.07179"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "0.6615"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "4.647"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").Value = "78.38"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "0.9989"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "26.166.79"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "11.91"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").Value = "0.000006813"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "1.984.49"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "4.421"
$ws.Range("E22").Value = "  +6.76%  "
$ws.Range("D23").Value = "8.790"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "5.311"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").Value = "140.30"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "1.516"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "15.33"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "1.809"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "105.67"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "0.08518"
$ws.Range("E30").Value = "  +3.00%  "
$ws.Range("D31").Value = "3.808"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("D32").Value = "3.674"
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("D33").Value = "0.04619"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("D34").Value = "2.663"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "1.006"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").Value = "0.6296"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").Value = "0.01618"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").Value = "1.961"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "100.09"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "0.3933"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").Value = "0.7536"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").Value = "5.048"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("D46").Value = "6.361"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "0.05359"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "55.07"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "31.02"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.702"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3510"
$ws.Range("E51").Value = "  +2.95%  "
